# Update "Top 50 Cryptocurrencies" sheet with refreshed live crypto data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

# Each entry: row, Name, Symbol, Current Price (USD), Market Capitalization, 24h Trading Volume, Price Change (24h %)
$rows = @(
    @(2, 'Bitcoin', 'btc', 98781, 1954198009844, 112390278976, 1.60471),
    @(3, 'Ethereum', 'eth', 3377.14, 406670065516, 57054771047, 8.46972),
    @(4, 'Tether', 'usdt', 1, 130821581010, 183012864316, -0.04268),
    @(5, 'Solana', 'sol', 259.48, 123188156523, 14777737320, 8.5897),
    @(6, 'BNB', 'bnb', 630.97, 92052602155, 2498055748, 3.47283),
    @(7, 'XRP', 'xrp', 1.41, 80370827021, 18422362798, 26.93999),
    @(8, 'Dogecoin', 'doge', 0.394378, 57906961638, 9730213490, 2.89369),
    @(9, 'USDC', 'usdc', 0.999152, 38292986605, 11690219708, -0.12898),
    @(10, 'Lido Staked Ether', 'steth', 3378.25, 33089387008, 141546953, 8.46048),
    @(11, 'Cardano', 'ada', 0.883571, 31609255024, 3758331072, 12.16342),
    @(12, 'TRON', 'trx', 0.199602, 17236971977, 1070924350, 1.23089),
    @(13, 'Avalanche', 'avax', 36.28, 14844705898, 1045261609, 7.19604),
    @(14, 'Shiba Inu', 'shib', 0.00002494, 14693630180, 1594894926, 3.59182),
    @(15, 'Wrapped stETH', 'wsteth', 4002.64, 14444484957, 169059682, 8.64603),
    @(16, 'Wrapped Bitcoin', 'wbtc', 98585, 14399440584, 826645767, 1.7654),
    @(17, 'Toncoin', 'ton', 5.54, 14119887111, 623059661, 3.12425),
    @(18, 'Sui', 'sui', 3.6, 10232927881, 2055446112, 0.63511),
    @(19, 'Bitcoin Cash', 'bch', 495.92, 9807988842, 1790730676, -4.17848),
    @(20, 'WETH', 'weth', 3377, 9674931817, 2223824701, 8.58658),
    @(21, 'Chainlink', 'link', 15.31, 9596106894, 1255002788, 5.58711),
    @(22, 'Pepe', 'pepe', 0.00002135, 8974266732, 6729179657, 10.04046),
    @(23, 'Polkadot', 'dot', 6.22, 8969229904, 840762108, 9.71723),
    @(24, 'Stellar', 'xlm', 0.286093, 8579890000, 2328650446, 19.85431),
    @(25, 'LEO Token', 'leo', 8.79, 8127976591, 3417151, 3.64561),
    @(26, 'NEAR Protocol', 'near', 5.8, 7067686644, 1004401919, 5.35271),
    @(27, 'Litecoin', 'ltc', 90.86, 6833649109, 1335029714, 3.91033),
    @(28, 'Aptos', 'apt', 12.09, 6472209222, 851267493, 3.9428),
    @(29, 'Wrapped eETH', 'weeth', 3559.51, 6194902543, 105782378, 8.69229),
    @(30, 'Uniswap', 'uni', 9.35, 5614203040, 870750030, 5.84469),
    @(31, 'Cronos', 'cro', 0.195065, 5341306291, 141273639, 11.78849),
    @(32, 'Hedera', 'hbar', 0.137038, 5238270224, 934498283, 10.09357),
    @(33, 'USDS', 'usds', 0.999299, 5230226962, 16080506, -0.32212),
    @(34, 'Internet Computer', 'icp', 9.66, 4584170683, 274022971, 7.12318),
    @(35, 'Ethereum Classic', 'etc', 28.06, 4199695234, 866959534, 5.17743),
    @(36, 'Bonk', 'bonk', 0.00005233, 3919437831, 1630787456, 0.41931),
    @(37, 'Kaspa', 'kas', 0.152056, 3832081221, 151141330, 0.92423),
    @(38, 'Render', 'render', 7.37, 3813426846, 430485438, 0.49782),
    @(39, 'POL (ex-MATIC)', 'pol', 0.467282, 3723235499, 497017761, 6.94723),
    @(40, 'Bittensor', 'tao', 502.34, 3722424094, 280112348, 3.21335),
    @(41, 'Ethena USDe', 'usde', 1.001, 3686884674, 224249364, -0.08049),
    @(42, 'WhiteBIT Coin', 'wbt', 24.77, 3575145979, 33131314, 2.74346),
    @(43, 'Dai', 'dai', 0.999429, 3439660368, 184363126, 0.14324),
    @(44, 'MANTRA', 'om', 3.82, 3439598122, 304755108, 5.80427),
    @(45, 'dogwifhat', 'wif', 3.38, 3377151669, 1281353338, 5.16064),
    @(46, 'Artificial Superintelligence Alliance', 'fet', 1.28, 3350991172, 481757661, 3.13522),
    @(47, 'Arbitrum', 'arb', 0.789976, 3236809423, 1676829198, 14.39195),
    @(48, 'Monero', 'xmr', 161.18, 2972719324, 83675039, -0.35097),
    @(49, 'Stacks', 'stx', 1.96, 2943890749, 352719354, 2.63933),
    @(50, 'Mantle', 'mnt', 0.838499, 2824529619, 185914513, 15.13845),
    @(51, 'Filecoin', 'fil', 4.69, 2815699050, 576199430, 7.35606)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws1.Cells.Item($rowNum, 1).Value = $r[1]
    $ws1.Cells.Item($rowNum, 2).Value = $r[2]
    $ws1.Cells.Item($rowNum, 3).Value = $r[3]
    $ws1.Cells.Item($rowNum, 4).Value = $r[4]
    $ws1.Cells.Item($rowNum, 5).Value = $r[5]
    $ws1.Cells.Item($rowNum, 6).Value = $r[6]
}

# Update "Top 5 by Market Cap" sheet (market caps mirror the top 5 rows above; names/order unchanged)
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws2.Cells.Item(2, 2).Value = 1954198009844
$ws2.Cells.Item(3, 2).Value = 406670065516
$ws2.Cells.Item(4, 2).Value = 130821581010
$ws2.Cells.Item(5, 2).Value = 123188156523
$ws2.Cells.Item(6, 2).Value = 92052602155

# Update "Summary" sheet
$ws3 = $wb.Worksheets.Item("Summary")
# Average Price must stay plain text (not auto-converted to a currency number)
$avgCell = $ws3.Cells.Item(2, 2)
$avgCell.NumberFormat = "@"
$avgCell.Value = '$4348.00'
$avgCell.NumberFormat = "General"
$avgCell.Style = "Normal"
$ws3.Cells.Item(3, 2).Value = 'XRP (26.94%)'
$ws3.Cells.Item(4, 2).Value = 'Bitcoin Cash (-4.18%)'
